$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.906.68"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.784.60"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.05"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.47"
$ws.Range("E6").Value = "  -3.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.558"
$ws.Range("E7").Value = "  -2.54%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.38"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.47"
$ws.Range("E13").Value = "  -3.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.57"
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("D15").Value = "3.225.49"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "2.794.88"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").Value = "51.846.44"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.50"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.15"
$ws.Range("E21").Value = "  -3.65%  "
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.16"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.01"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  -4.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.47"
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.163"
$ws.Range("E28").Value = "  +15.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.32"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0470"
$ws.Range("E31").Value = "  -3.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.10"
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.69"
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.74"
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("E36").Value = "  -5.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.84"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("E39").Value = "  -3.10%  "
$ws.Range("E40").Value = "  -4.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.63"
$ws.Range("E41").Value = "  +3.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.17"
$ws.Range("E44").Value = "  -4.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.81"
$ws.Range("E45").Value = "  -8.10%  "
$ws.Range("D46").Value = "2.084.93"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("E47").Value = "  -4.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.965"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("E51").Value = "  +31.82%  "
